# Insert a new data row at row 557 (pushes the existing row 557..603 down to 558..604)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(557).Insert()

# Populate the new row 557. Most descriptive columns (market/product metadata)
# are identical to the row that used to occupy 557 (now at 558); only the
# date, volume/price columns, unit, origin and $/Kg columns differ.
$ws.Cells.Item(557, 1).Value = 11
$ws.Cells.Item(557, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(557, 3).Value = "Bíobío"
$ws.Cells.Item(557, 4).Value = 45265
$ws.Cells.Item(557, 5).Value = 8
$ws.Cells.Item(557, 6).Value = "Fruta"
$ws.Cells.Item(557, 7).Value = 100101
$ws.Cells.Item(557, 8).Value = "Berries"
$ws.Cells.Item(557, 9).Value = 100112025
$ws.Cells.Item(557, 10).Value = "Frutilla"
$ws.Cells.Item(557, 11).Value = "Sin especificar"
$ws.Cells.Item(557, 12).Value = "Primera"
$ws.Cells.Item(557, 13).Value = 180
$ws.Cells.Item(557, 14).Value = 10000
$ws.Cells.Item(557, 15).Value = 11000
$ws.Cells.Item(557, 16).Value = 10556
$ws.Cells.Item(557, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(557, 18).Value = "Región del Maule"
$ws.Cells.Item(557, 19).Value = 1508
$ws.Cells.Item(557, 20).Value = 7
